$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = "[]"
$ws.Range("E13").Value = "[]"

$ws.Range("C16").Value = "[336]"
$ws.Range("E16").Value = "[336]"

$ws.Range("C18").Value = "[]"
$ws.Range("D18").Value = "[]"
$ws.Range("E18").Value = "[]"
$ws.Range("F18").Value = "[90]"

$ws.Range("C19").Value = "[585]"

$ws.Range("C20").Value = "[]"
$ws.Range("D20").Value = "[]"
$ws.Range("E20").Value = "[]"
$ws.Range("F20").Value = "[700]"

$ws.Range("C21").Value = "[]"
$ws.Range("D21").Value = "[]"
$ws.Range("E21").Value = "[]"
$ws.Range("F21").Value = "[260]"

$ws.Range("C22").Value = "[]"
$ws.Range("D22").Value = "[]"
$ws.Range("E22").Value = "[]"
$ws.Range("F22").Value = "[170]"

$ws.Range("C24").Value = "[]"
$ws.Range("D24").Value = "[]"
$ws.Range("E24").Value = "[]"
$ws.Range("F24").Value = "[340]"

$ws.Range("C26").Value = "[483]"

$ws.Range("C27").Value = "[375]"

$ws.Range("C28").Value = "[291]"

$ws.Range("C30").Value = "[150]"

$ws.Range("C31").Value = "[375]"
$ws.Range("D31").Value = "[300]"
$ws.Range("E31").Value = "[]"
$ws.Range("F31").Value = "[]"

$ws.Range("C33").Value = "[]"
$ws.Range("D33").Value = "[]"
$ws.Range("E33").Value = "[]"
$ws.Range("F33").Value = "[600]"

$ws.Range("C34").Value = "[468]"
$ws.Range("E34").Value = "[468]"

$ws.Range("C36").Value = "[213]"

$ws.Range("C37").Value = "[]"
$ws.Range("E37").Value = "[]"

$ws.Range("C38").Value = "[]"
$ws.Range("E38").Value = "[]"

$ws.Range("C39").Value = "[]"
$ws.Range("E39").Value = "[]"
